$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'35.303.00"
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(2, 5).Value = "'  -0.78%  "
$ws.Cells.Item(2, 5).ClearFormats()
$ws.Cells.Item(3, 4).Value = "'1.905.10"
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(3, 5).Value = "'  +0.88%  "
$ws.Cells.Item(3, 5).ClearFormats()
$ws.Cells.Item(4, 5).Value = "'  -0.28%  "
$ws.Cells.Item(4, 5).ClearFormats()
$ws.Cells.Item(5, 5).Value = "'  +9.47%  "
$ws.Cells.Item(5, 5).ClearFormats()
$ws.Cells.Item(6, 4).Value = "'245.89"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = "'  +0.43%  "
$ws.Cells.Item(6, 5).ClearFormats()
$ws.Cells.Item(7, 5).Value = "'  -0.22%  "
$ws.Cells.Item(7, 5).ClearFormats()
$ws.Cells.Item(8, 4).Value = "'41.60"
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).Value = "'  -3.35%  "
$ws.Cells.Item(8, 5).ClearFormats()
$ws.Cells.Item(9, 4).Value = "'0.349"
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).Value = "'  +4.36%  "
$ws.Cells.Item(9, 5).ClearFormats()
$ws.Cells.Item(10, 5).Value = "'  +11.78%  "
$ws.Cells.Item(10, 5).ClearFormats()
$ws.Cells.Item(11, 4).Value = "'0.0725"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value = "'  +2.64%  "
$ws.Cells.Item(11, 5).ClearFormats()
$ws.Cells.Item(12, 5).Value = "'  -0.01%  "
$ws.Cells.Item(12, 5).ClearFormats()
$ws.Cells.Item(13, 4).Value = "'2.181.19"
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).Value = "'  +1.08%  "
$ws.Cells.Item(13, 5).ClearFormats()
$ws.Cells.Item(14, 5).Value = "'  +0.44%  "
$ws.Cells.Item(14, 5).ClearFormats()
$ws.Cells.Item(15, 4).Value = "'0.712"
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).Value = "'  +3.13%  "
$ws.Cells.Item(15, 5).ClearFormats()
$ws.Cells.Item(16, 4).Value = "'1.905.83"
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(16, 5).Value = "'  +1.00%  "
$ws.Cells.Item(16, 5).ClearFormats()
$ws.Cells.Item(17, 4).Value = "'4.84"
$ws.Cells.Item(17, 4).ClearFormats()
$ws.Cells.Item(17, 5).Value = "'  +0.81%  "
$ws.Cells.Item(17, 5).ClearFormats()
$ws.Cells.Item(18, 4).Value = "'35.299.79"
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).Value = "'  -0.83%  "
$ws.Cells.Item(18, 5).ClearFormats()
$ws.Cells.Item(19, 4).Value = "'72.38"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = "'  +0.50%  "
$ws.Cells.Item(19, 5).ClearFormats()
$ws.Cells.Item(20, 4).Value = "'0.0₃0821"
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = "'  +1.41%  "
$ws.Cells.Item(20, 5).ClearFormats()
$ws.Cells.Item(21, 4).Value = "'241.15"
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).Value = "'  -1.16%  "
$ws.Cells.Item(21, 5).ClearFormats()
$ws.Cells.Item(22, 4).Value = "'12.68"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = "'  +1.89%  "
$ws.Cells.Item(22, 5).ClearFormats()
$ws.Cells.Item(23, 4).Value = "'4.82"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).Value = "'  -0.62%  "
$ws.Cells.Item(23, 5).ClearFormats()
$ws.Cells.Item(24, 5).Value = "'  -0.27%  "
$ws.Cells.Item(24, 5).ClearFormats()
$ws.Cells.Item(25, 4).Value = "'2.31"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = "'  +0.74%  "
$ws.Cells.Item(25, 5).ClearFormats()
$ws.Cells.Item(26, 4).Value = "'2.28"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).Value = "'  +8.64%  "
$ws.Cells.Item(26, 5).ClearFormats()
$ws.Cells.Item(27, 4).Value = "'169.09"
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value = "'  -1.05%  "
$ws.Cells.Item(27, 5).ClearFormats()
$ws.Cells.Item(28, 4).Value = "'8.59"
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).Value = "'  +2.32%  "
$ws.Cells.Item(28, 5).ClearFormats()
$ws.Cells.Item(29, 2).Value = "'Stellar"
$ws.Cells.Item(29, 2).ClearFormats()
$ws.Cells.Item(29, 3).Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(29, 3).ClearFormats()
$ws.Cells.Item(29, 4).Value = "'0.131"
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value = "'  +4.27%  "
$ws.Cells.Item(29, 5).ClearFormats()
$ws.Cells.Item(30, 2).Value = "'EthereumClassic"
$ws.Cells.Item(30, 2).ClearFormats()
$ws.Cells.Item(30, 3).Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(30, 3).ClearFormats()
$ws.Cells.Item(30, 4).Value = "'18.41"
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value = "'  +2.64%  "
$ws.Cells.Item(30, 5).ClearFormats()
$ws.Cells.Item(32, 5).Value = "'  +1.93%  "
$ws.Cells.Item(32, 5).ClearFormats()
$ws.Cells.Item(33, 4).Value = "'0.962"
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = "'  -1.85%  "
$ws.Cells.Item(33, 5).ClearFormats()
$ws.Cells.Item(34, 4).Value = "'0.0572"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = "'  +0.90%  "
$ws.Cells.Item(34, 5).ClearFormats()
$ws.Cells.Item(35, 5).Value = "'  -0.12%  "
$ws.Cells.Item(35, 5).ClearFormats()
$ws.Cells.Item(36, 4).Value = "'4.14"
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).Value = "'  +0.99%  "
$ws.Cells.Item(36, 5).ClearFormats()
$ws.Cells.Item(37, 5).Value = "'  +0.07%  "
$ws.Cells.Item(37, 5).ClearFormats()
$ws.Cells.Item(38, 4).Value = "'1.46"
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).Value = "'  +8.10%  "
$ws.Cells.Item(38, 5).ClearFormats()
$ws.Cells.Item(39, 5).Value = "'  -1.59%  "
$ws.Cells.Item(39, 5).ClearFormats()
$ws.Cells.Item(40, 5).Value = "'  +10.71%  "
$ws.Cells.Item(40, 5).ClearFormats()
$ws.Cells.Item(41, 4).Value = "'1.10"
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).Value = "'  -0.66%  "
$ws.Cells.Item(41, 5).ClearFormats()
$ws.Cells.Item(42, 5).Value = "'  +2.54%  "
$ws.Cells.Item(42, 5).ClearFormats()
$ws.Cells.Item(43, 4).Value = "'16.10"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = "'  +4.22%  "
$ws.Cells.Item(43, 5).ClearFormats()
$ws.Cells.Item(44, 4).Value = "'90.01"
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value = "'  -0.56%  "
$ws.Cells.Item(44, 5).ClearFormats()
$ws.Cells.Item(45, 4).Value = "'1.346.47"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = "'  -0.78%  "
$ws.Cells.Item(45, 5).ClearFormats()
$ws.Cells.Item(46, 4).Value = "'2.45"
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value = "'  +4.28%  "
$ws.Cells.Item(46, 5).ClearFormats()
$ws.Cells.Item(47, 4).Value = "'12.67"
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).Value = "'  -5.93%  "
$ws.Cells.Item(47, 5).ClearFormats()
$ws.Cells.Item(48, 2).Value = "'HuobiToken"
$ws.Cells.Item(48, 2).ClearFormats()
$ws.Cells.Item(48, 3).Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(48, 3).ClearFormats()
$ws.Cells.Item(48, 4).Value = "'2.42"
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).Value = "'  +0.10%  "
$ws.Cells.Item(48, 5).ClearFormats()
$ws.Cells.Item(49, 2).Value = "'MXToken"
$ws.Cells.Item(49, 2).ClearFormats()
$ws.Cells.Item(49, 3).Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(49, 3).ClearFormats()
$ws.Cells.Item(49, 4).Value = "'2.80"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value = "'  +1.48%  "
$ws.Cells.Item(49, 5).ClearFormats()
$ws.Cells.Item(50, 2).Value = "'MultiversX"
$ws.Cells.Item(50, 2).ClearFormats()
$ws.Cells.Item(50, 3).Value = "'https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Cells.Item(50, 3).ClearFormats()
$ws.Cells.Item(50, 4).Value = "'45.99"
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).Value = "'  +1.82%  "
$ws.Cells.Item(50, 5).ClearFormats()
$ws.Cells.Item(51, 5).Value = "'  -2.03%  "
$ws.Cells.Item(51, 5).ClearFormats()
